$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new ring row (row 43) to the worksheet data
$ws.Range("A43").Value = "Cloranty Ring"
$ws.Range("B43").Value = 3
$ws.Range("C43").Value = "tex_DaS_ClorantyRing.png"

# Grow the XML-mapped table to cover the new row, plus the trailing
# blank "insert" row Excel keeps under an XML list range
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C44"))

# Match the cursor position left behind after the XML list refresh
[void]$ws.Range("D44").Select()
